$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mirror the existing header formatting (bold font, border, centered alignment)
# from H1 onto the two new header cells I1 ("I0") and J1 ("IF").
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New numeric data columns I and J for rows 2-12
$values = @(
    @(9, 9),
    @(5, 6),
    @(6, 8),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(8, 9),
    @(8, 8),
    @(7, 8),
    @(6, 7),
    @(7, 7)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
